# usability_kevin.xlsx edit
# - Add a new observation row (row 13) in column C with a blank-ish
#   (two-space) remark, styled the same as the other "passed" cells
#   in that column (the built-in "Good" cell style).
# - Update the sheet's selection to C20:C21 (active cell anchored at
#   the top of the range, which is this engine's normalized behaviour
#   for a multi-cell .Select()).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New remark cell, reusing the "Good" style already applied to C2 so the
# appended shared string and style index line up with the target sheet.
$ws.Range("C13").Value = "  "
$ws.Range("C13").Style = $ws.Range("C2").Style

# Move/extend the selection to C20:C21.
$ws.Range("C20:C21").Select()
